$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns I (q18_important_features) and J (q39_contradictory_info)
# entirely, shifting nothing else - just dropping this unstructured data.
$ws.Range("I1:J46").EntireColumn.Delete()
